$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.168631076812744
$ws.Range("B1").Value = 1.625205516815186
$ws.Range("C1").Value = 1.286009907722473
$ws.Range("D1").Value = 1.886569857597351
$ws.Range("E1").Value = 3.143168449401855
